$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove existing hyperlinks on column F (worksheet-wide) so we can re-add clean ones ---
$ws.Range("F2").Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = "2026-01-16 06:31:19"
$ws.Range("B2").Value = "【募集】Python / Docker 日次データ スクレイピングシステム構築"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5469627"
$ws.Range("G2").Value = 248
$ws.Range("H2").Value = "🔥Python ◆スクレイピング"

# Row 3
$ws.Range("A3").Value = "2026-01-16 06:31:19"
$ws.Range("B3").Value = "【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5405023"
$ws.Range("G3").Value = 178
$ws.Range("H3").Value = "★bot ◆ツール"

# Row 4
$ws.Range("A4").Value = "2026-01-16 06:31:19"
$ws.Range("B4").Value = "【Windows/Wacom】署名画像から筆順解析図を作成する業務用アプリ開発"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5472804"
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = "◆開発 ◇アプリ"

# Row 5
$ws.Range("A5").Value = "2026-01-16 06:31:19"
$ws.Range("B5").Value = "署名画像から筆順を可視化するアプリ開発者募集"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5472080"
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = "◆開発 ◇アプリ"

# Row 6
$ws.Range("A6").Value = "2026-01-16 06:31:19"
$ws.Range("B6").Value = "スマホカラオケ予約Webアプリ開発のフリーランス募集(使用するのは個人の集まりで趣味で使う程度です)"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5472431"
$ws.Range("G6").Value = 88
$ws.Range("H6").Value = "◆開発 ◇アプリ"

# Row 7
$ws.Range("A7").Value = "2026-01-16 06:31:19"
$ws.Range("B7").Value = "初回 Webサーバ管理エンジニア"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5472544"
$ws.Range("G7").Value = 45
$ws.Range("H7").Value = "◇管理"

# Row 8
$ws.Range("A8").Value = "2026-01-16 06:31:19"
$ws.Range("B8").Value = "m.2 SSD基板の設計"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5472120"
$ws.Range("G8").Value = 25
$ws.Range("H8").ClearContents()

# --- Re-add hyperlinks for F2:F8 in order, then restore the shared Hyperlink cell style ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5469627")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5405023")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5472804")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5472080")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5472431")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5472544")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5472120")
$ws.Range("F2:F8").Style = "Hyperlink"

# --- Column B width: 41 -> 52 characters ---
$ws.Columns.Item(2).ColumnWidth = 51.16
